# Update the date strings held in column B of Sheet1 and the current
# selection, matching the authored change (dates shifted from the 3rd to
# the 6th of August 2023, and the selected range updated to B16:B17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use a leading apostrophe via Formula so the values are stored as plain
# text (matching the existing quote-prefixed text styles already applied
# to these cells) instead of being re-interpreted/reformatted as dates.

# These cells store "DD/MM/YYYY" style text -> "03/08/2023" becomes "06/08/2023"
$cellsDMY = @("B1","B2","B3","B4","B5","B6","B7","B8","B10","B11","B16","B17")
foreach ($addr in $cellsDMY) {
    $ws.Range($addr).Formula = "'06/08/2023"
}

# B12 stores "YYYY/MM/DD" style text -> "2023/08/03" becomes "2023/08/06"
$ws.Range("B12").Formula = "'2023/08/06"

# B14 stores "MM/DD/YYYY" style text -> "08/03/2023" becomes "08/06/2023"
$ws.Range("B14").Formula = "'08/06/2023"

# Update the active selection on the sheet to match the new state:
# B16:B17 selected, with B16 as the active cell.
$ws.Range("B16:B17").Select()
